# The commit bumps the "Förändrad" (Changed) date in column C by one day
# (from serial 46061 = 2026-02-08 to serial 46062 = 2026-02-09) for every
# data row in the sheet (rows 2 through 490).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C490").Value2 = 46062
